$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.084468483924866
$ws.Range("B1").Value = 1.263489723205566
$ws.Range("C1").Value = 1.634562849998474
$ws.Range("D1").Value = 3.64272141456604
$ws.Range("E1").Value = 3.789305925369263
